{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// The footer block to remove is the trailing run of paragraphs that\n// appears right after \"LOQ4203: Sistemas Produtivos I (Requisito fraco)\":\n//   1) an empty spacer paragraph\n//   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3) the \"\u00a9 2020 . Contact: ...\" copyright paragraph\n// A further empty paragraph (and the page-break paragraph after it) stay.\nconst items = paragraphs.items;\nlet reqIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === \"LOQ4203: Sistemas Produtivos I (Requisito fraco)\") {\n    reqIndex = i;\n    break;\n  }\n}\n\nif (reqIndex !== -1) {\n  const verText = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\n  const copyrightPrefix = \"\u00a9 2020\";\n\n  // Expect: reqIndex+1 empty, reqIndex+2 \"Ver no Jupiter...\", reqIndex+3 copyright.\n  if (\n    items[reqIndex + 1] &&\n    items[reqIndex + 1].text === \"\" &&\n    items[reqIndex + 2] &&\n    items[reqIndex + 2].text === verText &&\n    items[reqIndex + 3] &&\n    items[reqIndex + 3].text.indexOf(copyrightPrefix) === 0\n  ) {\n    // Delete in reverse order so earlier indices stay valid.\n    items[reqIndex + 3].delete();\n    items[reqIndex + 2].delete();\n    items[reqIndex + 1].delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"LOQ4203: Sistemas Produtivos I (Requisito fraco)\" paragraph,\n# then remove the trailing footer block that follows it:\n#   1) an empty spacer paragraph\n#   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3) the copyright paragraph (\"(c) 2020 . Contact: luizeleno@usp.br. ...\")\n# One more empty paragraph (and the page-break paragraph after it) stay untouched.\n$reqIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $txt = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($txt -eq \"LOQ4203: Sistemas Produtivos I (Requisito fraco)\") {\n        $reqIndex = $i\n        break\n    }\n}\n\nif ($reqIndex -ne -1) {\n    $verText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n\n    $t1 = $d.Paragraphs.Item($reqIndex + 1).Range.Text.TrimEnd([char]13, [char]7)\n    $t2 = $d.Paragraphs.Item($reqIndex + 2).Range.Text.TrimEnd([char]13, [char]7)\n    $t3 = $d.Paragraphs.Item($reqIndex + 3).Range.Text.TrimEnd([char]13, [char]7)\n\n    if ($t1 -eq \"\" -and $t2 -eq $verText -and $t3.Contains(\"Contact: luizeleno@usp.br\")) {\n        # Delete from the last paragraph back to the first so earlier indices\n        # remain valid while the later ones are removed.\n        $d.Paragraphs.Item($reqIndex + 3).Range.Delete()\n        $d.Paragraphs.Item($reqIndex + 2).Range.Delete()\n        $d.Paragraphs.Item($reqIndex + 1).Range.Delete()\n    }\n}\n"}
